# Consolidate Cargo Type & Vehicle Type on the QSfV (Quantization Size for
# Vehicles) sheet: instead of separate "passengers"/"freight" columns, each
# vehicle type gets two rows - one for "passenger <type>" and one for
# "freight <type>" - sharing a single "vehicles" quantity column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QSfV")

# Capture the original (pre-edit) row 2-7 data before we start overwriting it.
$types  = @()
$pass   = @()
$freight = @()
for ($r = 2; $r -le 7; $r++) {
    $types   += $ws.Cells.Item($r, 1).Value2
    $pass    += $ws.Cells.Item($r, 2).Value2
    $freight += $ws.Cells.Item($r, 3).Value2
}

# Clear the old "freight" header; column C (previously "freight") goes away
# entirely once the data below has been migrated into new rows.
$ws.Range("C1").ClearContents()

# Rows 2-7 become the "passenger <type>" rows, keeping the old passenger
# quantities in column B.
for ($i = 0; $i -lt 6; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = "passenger " + $types[$i]
    $ws.Cells.Item($r, 2).Value = $pass[$i]
}

# Rows 8-13 become the new "freight <type>" rows, using the quantities that
# used to live in column C.
for ($i = 0; $i -lt 6; $i++) {
    $r = 8 + $i
    $ws.Cells.Item($r, 1).Value = "freight " + $types[$i]
    $ws.Cells.Item($r, 2).Value = $freight[$i]
}

# Header row: column B now just holds "vehicles" (a single combined count
# covering both the old "passengers" and "freight" columns).
$ws.Range("B1").Value = "vehicles"

# Drop the now-empty column C entirely.
$ws.Columns.Item(3).Delete()
